$d = $word.ActiveDocument

# --- Locate the anchor paragraph: the list item ending in "...How spilit."
#     which currently also carries the (hidden) _GoBack bookmark at its end.
$howSpilitIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*How spilit.*") {
        $howSpilitIndex = $i
        break
    }
}

# --- Locate the existing "Start with predetermined list of charities"
#     paragraph (the one immediately following the anchor paragraph).
$startWithIndex = -1
for ($i = $howSpilitIndex; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Start with predetermined list of charities*") {
        $startWithIndex = $i
        break
    }
}

# 1) Remove the _GoBack bookmark from the anchor paragraph; it will be
#    re-homed to its own standalone paragraph below.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2) Duplicate the "Start with predetermined list of charities" paragraph:
#    insert a brand-new paragraph directly before the existing one (it
#    inherits the ListParagraph/numbering formatting of that paragraph),
#    then give it the same text. This new paragraph becomes the item right
#    after the "...How spilit." paragraph, matching the target diff.
$existingStartWith = $d.Paragraphs.Item($startWithIndex)
$existingStartWith.Range.InsertParagraphBefore()
$newStartWith = $d.Paragraphs.Item($startWithIndex)
$newStartWith.Range.Text = "Start with predetermined list of charities"

# 3) The paragraph that used to hold "Start with predetermined..." has now
#    shifted one slot later - that is the now-redundant duplicate. Directly
#    after it sits the paragraph holding only <w:ind w:left="360"/>, and
#    after that a bare empty paragraph. Delete the duplicate item and the
#    ind=360 paragraph, leaving the bare empty paragraph in place so it can
#    become the new home for the _GoBack bookmark.
$dupIndex = $startWithIndex + 1
$d.Paragraphs.Item($dupIndex).Range.Delete()
$d.Paragraphs.Item($dupIndex).Range.Delete()

# 4) Re-insert the _GoBack bookmark into that now-bare empty paragraph.
$bookmarkRange = $d.Paragraphs.Item($dupIndex).Range
$bookmarkRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
